# Actualizar 05-09-2020 04-16-43
# Updates the HOSPITALES sheet: renames a duplicated hospital/clinic name,
# extends the HOSPITALES_HN table, and appends new hospital/clinic records
# for Yoro and Lempira departments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix duplicated "San Lucas" entry: the clinic in Gracias (Lempira) and
# the clinic in Catacamas (Olancho) should both read "Clinica San Lucas".
$ws.Range("U47").Value = 'Clínica San Lucas'
$ws.Range("U102").Value = 'Clínica San Lucas'

# --- Append new rows (122-129) with new health-center records ---

# Row 122: Yoro (department, with trailing space variant) - Centro de Salud San Antonio Sulaco Yoro
$ws.Range("G122").Value = 'Yoro '
$ws.Range("S122").Value = 'centro de salud'
$ws.Range("U122").Value = 'Centro de Salud San Antonio Sulaco Yoro'
$ws.Range("V122").Value = 14.985823
$ws.Range("W122").Value = -87.294527

# Row 123: Yoro / Yorito - Centro de Salud Yorito
$ws.Range("G123").Value = 'Yoro'
$ws.Range("K123").Value = 'Yorito'
$ws.Range("S123").Value = 'centro de salud'
$ws.Range("U123").Value = 'Centro de Salud Yorito'
$ws.Range("V123").Value = 15.065429
$ws.Range("W123").Value = -87.277243

# Row 124: Yoro - Centro de Salud La Sabana
$ws.Range("G124").Value = 'Yoro'
$ws.Range("S124").Value = 'centro de salud'
$ws.Range("U124").Value = 'Centro de Salud La Sabana'
$ws.Range("V124").Value = 15.108616
$ws.Range("W124").Value = -87.280531

# Row 125: Yoro / Yoro - Clinica Medica Dra. Chavez
$ws.Range("G125").Value = 'Yoro'
$ws.Range("K125").Value = 'Yoro'
$ws.Range("S125").Value = 'clínica'
$ws.Range("U125").Value = 'Clínica Médica Dra. Chávez'
$ws.Range("V125").Value = 15.137856
$ws.Range("W125").Value = -87.125667

# Row 126: Yoro / Yoro - Hospital Manuel de Jesus Subirana
$ws.Range("G126").Value = 'Yoro'
$ws.Range("K126").Value = 'Yoro'
$ws.Range("S126").Value = 'hospital'
$ws.Range("U126").Value = 'Hospital Manuel de Jesus Subirana'
$ws.Range("V126").Value = 15.137337
$ws.Range("W126").Value = -87.133142

# Row 127: Yoro / El Rosario - ACTS Clinic
$ws.Range("G127").Value = 'Yoro'
$ws.Range("K127").Value = 'El Rosario'
$ws.Range("S127").Value = 'clínica'
$ws.Range("U127").Value = 'ACTS Clinic'
$ws.Range("V127").Value = 15.276194
$ws.Range("W127").Value = -87.324968

# Row 128: Lempira / Lepaera - Centro de Salud Lepaera
$ws.Range("G128").Value = 'Lempira'
$ws.Range("K128").Value = 'Lepaera'
$ws.Range("S128").Value = 'centro de salud'
$ws.Range("U128").Value = 'Centro de Salud Lepaera'
$ws.Range("V128").Value = 14.78034
$ws.Range("W128").Value = -88.588607

# Row 129: Cesamo de Santa Cruz (category only, no department/municipality set)
$ws.Range("S129").Value = 'centro de salud'
$ws.Range("U129").Value = 'Cesamo de Santa Cruz '
$ws.Range("V129").Value = 14.329539
$ws.Range("W129").Value = -88.522635

# --- Resize the HOSPITALES_HN table / autofilter to cover the new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:W131"))

# --- Update the hidden _FilterDatabase defined name to match ---
$fdb = $wb.Names.Item("_xlnm._FilterDatabase")
$fdb.RefersTo = "=HOSPITALES!`$A`$1:`$W`$131"

# --- Restore the active selection near the newly added data ---
$ws.Range("U131").Select() | Out-Null
